# "Worked on Adding codewords in existing codewordset"
# Replace the existing word list in column A (A1:A10) with the new
# "Testcodeword1" .. "Testcodeword10" codeword set, then resize the
# column to fit the new (longer) text and leave the selection on the
# next free cells below the data (A11:A14), matching where a user would
# continue typing more codewords.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "Testcodeword1",
    "Testcodeword2",
    "Testcodeword3",
    "Testcodeword4",
    "Testcodeword5",
    "Testcodeword6",
    "Testcodeword7",
    "Testcodeword8",
    "Testcodeword9",
    "Testcodeword10"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $words[$i]
}

# Column A needs to widen to fit the longer "Testcodeword.." strings.
$ws.Columns("A:A").AutoFit() | Out-Null

# Leave the selection below the new data, ready for more entries.
$ws.Range("A11:A14").Select() | Out-Null
